$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "I am the life of the party.",
        "ques_type": 2,
        "options": [
            "Very accurate.",
            "Moderately accurate.",
            "Neither inaccurate nor accurate.",
            "Moderately inaccurate.",
            "Very inaccurate."
        ],
        "score": "Very accurate."
    },
    {
        "title": "I feel little concern for others.",
        "ques_type": 2,
        "options": [
            "Very accurate.",
            "Moderately accurate.",
            "Neither inaccurate nor accurate.",
            "Moderately inaccurate.",
            "Very inaccurate."
        ],
        "score": "Very accurate."
    },
    {
        "title": "I am always prepared.",
        "ques_type": 2,
        "options": [
            "Very accurate.",
            "Moderately accurate.",
            "Neither inaccurate nor accurate.",
            "Moderately inaccurate.",
            "Very inaccurate."
        ],
        "score": "Very accurate."
    },
    {
        "title": "I get stressed out easily.",
        "ques_type": 2,
        "options": [
            "Very accurate.",
            "Moderately accurate.",
            "Neither inaccurate nor accurate.",
            "Moderately inaccurate.",
            "Very inaccurate."
        ],
        "score": "Very accurate."
    }
]
'@

# Remove the old second row entirely so the used range collapses back to A1
$ws.Range("A2").ClearContents()

# A1 previously held a bold/bordered "0" - strip that formatting before
# writing the new text so the cell reverts to the default (unstyled) xf
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text

# Undo the implicit autofit-row-height bump that comes from dropping a
# multi-line string into the cell, so row 1 keeps the sheet's default height
$ws.Rows(1).AutoFit()
